$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.577.63"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.70%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.524.12"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -1.73%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.05%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'304.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.75%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'96.57"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -1.13%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.62%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  +0.11%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.537"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -2.10%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'36.41"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -0.54%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  -0.17%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'7.49"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -2.10%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  -1.39%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'2.910.43"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -1.83%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'2.508.36"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -3.82%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'15.05"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +4.59%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'  -2.70%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'42.593.45"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.69%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'12.92"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +0.26%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  -2.46%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'  -2.69%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  -1.25%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'251.08"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -1.39%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'2.91"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -1.49%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'2.01"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -5.13%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'27.08"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -5.72%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  +0.18%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  +10.33%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  +0.13%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'38.23"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +0.65%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  -1.32%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'155.16"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +0.12%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'3.32"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -2.51%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'0.0789"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -1.80%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("B35").Value = "'Celestia"
$ws.Range("B35").Style = "Normal"
$ws.Range("C35").Value = "'https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("C35").Style = "Normal"
$ws.Range("D35").Value = "'18.72"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +2.50%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  -4.37%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("B37").Value = "'WEMIXToken"
$ws.Range("B37").Style = "Normal"
$ws.Range("C37").Value = "'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("C37").Style = "Normal"
$ws.Range("D37").Value = "'2.63"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -4.44%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.115"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +1.01%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("B39").Value = "'Stellar"
$ws.Range("B39").Style = "Normal"
$ws.Range("C39").Value = "'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("C39").Style = "Normal"
$ws.Range("D39").Value = "'0.119"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -0.73%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("B40").Value = "'EnergySwap"
$ws.Range("B40").Style = "Normal"
$ws.Range("C40").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("C40").Style = "Normal"
$ws.Range("D40").Value = "'23.96"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +3.26%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'3.39"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -1.01%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("B42").Value = "'RenderToken"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = "'3.84"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -1.03%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("B43").Value = "'FirstDigitalUSD"
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = "'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = "'0.998"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -0.07%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("B44").Value = "'ApeXProtocol"
$ws.Range("B44").Style = "Normal"
$ws.Range("C44").Value = "'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("C44").Style = "Normal"
$ws.Range("D44").Value = "'2.02"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -2.74%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  -3.31%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'2.026.91"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -2.06%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'85.12"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -0.26%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'8.96"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -2.62%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'2.765.82"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -2.00%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  -0.98%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'101.74"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -4.53%  "
$ws.Range("E51").Style = "Normal"
